# Update the date and all the two-digit multiplication problems/answers
# in the document to match the new "output generated at aa3dc9e" content.

$d = $word.ActiveDocument

# Mapping of old text -> new text (all values are unique within the document,
# so a straightforward Find/Replace for each pair is safe regardless of order).
$replacements = @(
    @("2024-01-03 Wednesday", "2024-01-04 Thursday"),
    @("21×36=756",  "78×60=4680"),
    @("27×50=1350", "30×84=2520"),
    @("90×29=2610", "84×55=4620"),
    @("30×57=1710", "96×39=3744"),
    @("97×49=4753", "38×39=1482"),
    @("69×29=2001", "86×95=8170"),
    @("54×46=2484", "28×30=840"),
    @("40×23=920",  "67×29=1943"),
    @("32×37=1184", "91×32=2912"),
    @("89×15=1335", "72×72=5184"),
    @("55×59=3245", "67×13=871"),
    @("77×88=6776", "66×94=6204"),
    @("82×30=2460", "47×84=3948"),
    @("69×14=966",  "31×22=682"),
    @("30×38=1140", "44×20=880"),
    @("60×91=5460", "81×97=7857"),
    @("33×28=924",  "98×36=3528"),
    @("16×58=928",  "81×88=7128"),
    @("27×26=702",  "59×36=2124"),
    @("32×68=2176", "75×87=6525"),
    @("92×55=5060", "78×46=3588"),
    @("72×50=3600", "24×75=1800"),
    @("77×14=1078", "57×30=1710"),
    @("67×64=4288", "16×27=432"),
    @("57×72=4104", "51×95=4845")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
